# Updated CPlusPlusCodeParser to handle SRC.<type>_decl tags.
#
# Mirrors the target diff against ABB.SrcML.Data/SrcMLTagParsing.xlsx:
#   - Sheet2!F5/F6 ("DeclarationStatement" rows): TODO/ParseExpressionStatementElement
#     -> "ParseDeclarationStatementElement"
#   - Sheet2!G6 comment expanded to mention ParsePropertyDeclarationElement
#   - Sheet2!F16 ("Foreach"/C#): TODO -> "ParseForeachElement"
#   - Sheet2!F35,F37,F39,F41,F42,F43 (C++ *Declaration / access-specifier rows):
#     TODO -> "ParseTypeElement"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Order matters for shared-string append order (new strings are appended to
# xl/sharedStrings.xml in first-use order), so touch ParseForeachElement,
# then ParseDeclarationStatementElement, then the expanded comment - same
# order as the new <si> entries appear in the target sharedStrings.xml.
$ws.Range("F16").Value = "ParseForeachElement"

$ws.Range("F5").Value = "ParseDeclarationStatementElement"
$ws.Range("F6").Value = "ParseDeclarationStatementElement"
$ws.Range("G6").Value = "Special handling for properties. Handled in ParsePropertyDeclarationElement"

$ws.Range("F35").Value = "ParseTypeElement"
$ws.Range("F37").Value = "ParseTypeElement"
$ws.Range("F39").Value = "ParseTypeElement"
$ws.Range("F41").Value = "ParseTypeElement"
$ws.Range("F42").Value = "ParseTypeElement"
$ws.Range("F43").Value = "ParseTypeElement"

# Best-effort: reproduce the scrolled/selected view state from the diff
# (frozen header pane scrolled so row 26 is the first visible data row,
# with F39 as the active cell).
$win = $excel.ActiveWindow
$ws.Range("A26").Select()
$win.ScrollRow = 26
$win.ScrollColumn = 1
$ws.Range("F39").Select()
